$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.255343993491408
$ws.Range("C2").Value = 0.2792743365243666
$ws.Range("D2").Value = 0.03037062535724289
$ws.Range("F2").Value = 1.086207185264144
$ws.Range("G2").Value = 0.9535733052840811
$ws.Range("H2").Value = 0.9431533495270514
$ws.Range("L2").Value = 0.2556585119273294
$ws.Range("M2").Value = 0.2832470116721026
$ws.Range("N2").Value = 1.505909118248321
$ws.Range("B3").Value = 1.156498875534112
$ws.Range("C3").Value = 0.2622586660347963
$ws.Range("D3").Value = 0.03040761311481965
$ws.Range("F3").Value = 1.069595557748599
$ws.Range("G3").Value = 0.9368146487658606
$ws.Range("H3").Value = 0.9409780274821742
$ws.Range("L3").Value = 0.2538986823195444
$ws.Range("M3").Value = 0.2687290893721581
$ws.Range("N3").Value = 1.525592623028956
$ws.Range("B4").Value = 1.096294837006042
$ws.Range("C4").Value = 0.2517087979723271
$ws.Range("D4").Value = 0.03043435331789723
$ws.Range("F4").Value = 1.060163038330217
$ws.Range("G4").Value = 0.9272567906128018
$ws.Range("H4").Value = 0.9401913594203108
$ws.Range("L4").Value = 0.2529435610165365
$ws.Range("M4").Value = 0.2599543572554595
$ws.Range("N4").Value = 1.538279445516712
$ws.Range("B5").Value = 1.071884028447982
$ws.Range("C5").Value = 0.2473839451689344
$ws.Range("D5").Value = 0.0304462682179949
$ws.Range("F5").Value = 1.056511545755995
$ws.Range("G5").Value = 0.9235452452521855
$ws.Range("H5").Value = 0.9400086318976122
$ws.Range("L5").Value = 0.2525858957110145
$ws.Range("M5").Value = 0.2564136629528591
$ws.Range("N5").Value = 1.543600416714283
$ws.Range("B6").Value = 1.067838059809617
$ws.Range("C6").Value = 0.2466642546530124
$ws.Range("D6").Value = 0.03044830830861045
$ws.Range("F6").Value = 1.055916818653714
$ws.Range("G6").Value = 0.9229399996122254
$ws.Range("H6").Value = 0.93998660985595
$ws.Range("L6").Value = 0.2525284119600784
$ws.Range("M6").Value = 0.2558278534846963
$ws.Range("N6").Value = 1.544493074733689
$ws.Range("B7").Value = 1.095965126584417
$ws.Range("C7").Value = 0.2516505755360754
$ws.Range("D7").Value = 0.03043450987835072
$ws.Range("F7").Value = 1.060113015002273
$ws.Range("G7").Value = 0.9272059939460036
$ws.Range("H7").Value = 0.9401883372312057
$ws.Range("L7").Value = 0.2529386096322312
$ws.Range("M7").Value = 0.2599064641151472
$ws.Range("N7").Value = 1.538350594862729
$ws.Range("B8").Value = 1.221161209747265
$ws.Range("C8").Value = 0.2734285338881079
$ws.Range("D8").Value = 0.0303825452009292
$ws.Range("F8").Value = 1.08031994446803
$ws.Range("G8").Value = 0.9476425505016977
$ws.Range("H8").Value = 0.942289230773099
$ws.Range("L8").Value = 0.255025699867538
$ws.Range("M8").Value = 0.2782123349818235
$ws.Range("N8").Value = 1.512571146017941
$ws.Range("B9").Value = 1.470537372940953
$ws.Range("C9").Value = 0.3153263499823424
$ws.Range("D9").Value = 0.03031241879390834
$ws.Range("F9").Value = 1.126060883041248
$ws.Range("G9").Value = 0.9935625832988535
$ws.Range("H9").Value = 0.9507757432033372
$ws.Range("L9").Value = 0.2601134825331357
$ws.Range("M9").Value = 0.3152158999925376
$ws.Range("N9").Value = 1.466793162016196
$ws.Range("B10").Value = 1.656133259617434
$ws.Range("C10").Value = 0.3456212717600806
$ws.Range("D10").Value = 0.03028002392071727
$ws.Range("F10").Value = 1.163439847176505
$ws.Range("G10").Value = 1.030917199432707
$ws.Range("H10").Value = 0.9596899312896596
$ws.Range("L10").Value = 0.2644586739366872
$ws.Range("M10").Value = 0.3430802309991634
$ws.Range("N10").Value = 1.436078705242718
$ws.Range("B11").Value = 1.741088333008918
$ws.Range("C11").Value = 0.359298759115859
$ws.Range("D11").Value = 0.0302693924154589
$ws.Range("F11").Value = 1.181273884209745
$ws.Range("G11").Value = 1.048708303464139
$ws.Range("H11").Value = 0.9643308532807282
$ws.Range("L11").Value = 0.2665674574186596
$ws.Range("M11").Value = 0.3559046344989483
$ws.Range("N11").Value = 1.422741283040908
$ws.Range("B12").Value = 1.773334386584565
$ws.Range("C12").Value = 0.3644631866777104
$ws.Range("D12").Value = 0.03026595309036395
$ws.Range("F12").Value = 1.188147259696436
$ws.Range("G12").Value = 1.055561007239646
$ws.Range("H12").Value = 0.9661727660405859
$ws.Range("L12").Value = 0.2673850006329275
$ws.Range("M12").Value = 0.3607823112222661
$ws.Range("N12").Value = 1.417782208430914
$ws.Range("B13").Value = 1.766386272892476
$ws.Range("C13").Value = 0.3633515994147558
$ws.Range("D13").Value = 0.03026666777865827
$ws.Range("F13").Value = 1.186661607524371
$ws.Range("G13").Value = 1.054080000739106
$ws.Range("H13").Value = 0.9657723150494064
$ws.Range("L13").Value = 0.2672080834915249
$ws.Range("M13").Value = 0.3597308672005752
$ws.Range("N13").Value = 1.418846159672348
$ws.Range("B14").Value = 1.743739725633361
$ws.Range("C14").Value = 0.3597239391546339
$ws.Range("D14").Value = 0.03026909772473374
$ws.Range("F14").Value = 1.181836952053146
$ws.Range("G14").Value = 1.049269758137029
$ws.Range("H14").Value = 0.9644806936373982
$ws.Range("L14").Value = 0.2666343365932846
$ws.Range("M14").Value = 0.3563054961722116
$ws.Range("N14").Value = 1.422331461425853
$ws.Range("B15").Value = 1.729877873248427
$ws.Range("C15").Value = 0.3574999465772919
$ws.Range("D15").Value = 0.0302706624221436
$ws.Range("F15").Value = 1.178897359430792
$ws.Range("G15").Value = 1.046338423952506
$ws.Range("H15").Value = 0.9637005498200324
$ws.Range("L15").Value = 0.2662853732300903
$ws.Range("M15").Value = 0.3542101371746895
$ws.Range("N15").Value = 1.424478234999494
$ws.Range("B16").Value = 1.650591815103837
$ws.Range("C16").Value = 0.3447253256847205
$ws.Range("D16").Value = 0.0302808009972253
$ws.Range("F16").Value = 1.162291110190083
$ws.Range("G16").Value = 1.029770641455116
$ws.Range("H16").Value = 0.9593984460956335
$ws.Range("L16").Value = 0.2643235181509596
$ws.Range("M16").Value = 0.342245116364829
$ws.Range("N16").Value = 1.436963130368402
$ws.Range("B17").Value = 1.602087031405688
$ws.Range("C17").Value = 0.3368619100044157
$ws.Range("D17").Value = 0.03028806939443385
$ws.Range("F17").Value = 1.152316789738123
$ws.Range("G17").Value = 1.019811878610284
$ws.Range("H17").Value = 0.9569094676756436
$ws.Range("L17").Value = 0.2631538204032751
$ws.Range("M17").Value = 0.3349430400722326
$ws.Range("N17").Value = 1.44478495230365
$ws.Range("B18").Value = 1.574237944253071
$ws.Range("C18").Value = 0.3323293093892232
$ws.Range("D18").Value = 0.03029263658092418
$ws.Range("F18").Value = 1.146657924580524
$ws.Range("G18").Value = 1.0141589532588
$ws.Range("H18").Value = 0.9555329970589241
$ws.Range("L18").Value = 0.2624934787305335
$ws.Range("M18").Value = 0.3307570851028956
$ws.Range("N18").Value = 1.449343587068709
$ws.Range("B19").Value = 1.56481724711864
$ws.Range("C19").Value = 0.3307929706941479
$ws.Range("D19").Value = 0.030294249467957
$ws.Range("F19").Value = 1.144755324508395
$ws.Range("G19").Value = 1.012257842496155
$ws.Range("H19").Value = 0.9550764075675033
$ws.Range("L19").Value = 0.26227203479975
$ws.Range("M19").Value = 0.3293421995216761
$ws.Range("N19").Value = 1.450897313788493
$ws.Range("B20").Value = 1.607245319206129
$ws.Range("C20").Value = 0.3376999954616338
$ws.Range("D20").Value = 0.03028725567628499
$ws.Range("F20").Value = 1.153370485498655
$ws.Range("G20").Value = 1.020864228398978
$ws.Range("H20").Value = 0.9571687166964864
$ws.Range("L20").Value = 0.2632770495481935
$ws.Range("M20").Value = 0.3357189092735169
$ws.Range("N20").Value = 1.443946122456726
$ws.Range("B21").Value = 1.750389521310467
$ws.Range("C21").Value = 0.3607898764550157
$ws.Range("D21").Value = 0.03026836809857159
$ws.Range("F21").Value = 1.183250808111225
$ws.Range("G21").Value = 1.050679499246627
$ws.Range("H21").Value = 0.9648577788861417
$ws.Range("L21").Value = 0.2668023445702232
$ws.Range("M21").Value = 0.3573110314669421
$ws.Range("N21").Value = 1.421305258086866
$ws.Range("B22").Value = 1.844381948618775
$ws.Range("C22").Value = 0.3757933623300573
$ws.Range("D22").Value = 0.0302594415147972
$ws.Range("F22").Value = 1.203479150442817
$ws.Range("G22").Value = 1.07083967643959
$ws.Range("H22").Value = 0.970375656068228
$ws.Range("L22").Value = 0.269217029208221
$ws.Range("M22").Value = 0.3715471779076722
$ws.Range("N22").Value = 1.407041696290629
$ws.Range("B23").Value = 1.794176314444655
$ws.Range("C23").Value = 0.3677936910548567
$ws.Range("D23").Value = 0.03026389427145659
$ws.Range("F23").Value = 1.192618658076981
$ws.Range("G23").Value = 1.060017856833042
$ws.Range("H23").Value = 0.967385499773286
$ws.Range("L23").Value = 0.2679181399651043
$ws.Range("M23").Value = 0.3639377048627566
$ws.Range("N23").Value = 1.414605525931535
$ws.Range("B24").Value = 1.604913142147723
$ws.Range("C24").Value = 0.3373211338442843
$ws.Range("D24").Value = 0.03028762234736604
$ws.Range("F24").Value = 1.152893874495291
$ws.Range("G24").Value = 1.020388235284059
$ws.Range("H24").Value = 0.9570513405328711
$ws.Range("L24").Value = 0.2632212998586994
$ws.Range("M24").Value = 0.3353681010963498
$ws.Range("N24").Value = 1.444325165039284
$ws.Range("B25").Value = 1.402658305779141
$ws.Range("C25").Value = 0.3040778817314731
$ws.Range("D25").Value = 0.03032801400599006
$ws.Range("F25").Value = 1.113027423437018
$ws.Range("G25").Value = 0.9805087363982921
$ws.Range("H25").Value = 0.9480106765633423
$ws.Range("L25").Value = 0.258630519019114
$ws.Range("M25").Value = 0.3050866963113776
$ws.Range("N25").Value = 1.478665332361262
